$wb = $excel.ActiveWorkbook

# Sheet ALC, row 62  (@@ -3736,25 +3736,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 125006936
$ws.Range("I62").Value = 55565236
$ws.Range("J62").Value = 250002000
$ws.Range("K62").Value = 55565236
$ws.Range("L62").Value = 250002000
$ws.Range("M62").Value = -55564612
$ws.Range("N62").Value = -250003248

# Sheet ALC, row 65  (@@ -3886,25 +3886,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 125006936
$ws.Range("I65").Value = 55565236
$ws.Range("J65").Value = 250002000
$ws.Range("K65").Value = 277826180
$ws.Range("L65").Value = 1250010000
$ws.Range("M65").Value = -277823060
$ws.Range("N65").Value = -1250016240

# Sheet ARM, row 32  (@@ -9283,25 +9283,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41388.89
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 41388.89
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 41388.89
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -41962.89

# Sheet ARM, row 37  (@@ -9525,25 +9522,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 9387.223
$ws.Range("J37").Value = 14674.777
$ws.Range("L37").Value = 14674.777
$ws.Range("N37").Value = -15220.777

# Sheet ARM, row 61  (@@ -10695,25 +10692,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2372769
$ws.Range("I61").Value = 1386163.2
$ws.Range("J61").Value = 6538438
$ws.Range("K61").Value = 1386163.2
$ws.Range("L61").Value = 6538438
$ws.Range("M61").Value = -1385951.2
$ws.Range("N61").Value = -6538862

# Sheet ARM, row 74  (@@ -11347,25 +11344,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 11853663
$ws.Range("I74").Value = 786.23334
$ws.Range("J74").Value = 35559416
$ws.Range("K74").Value = 786.23334
$ws.Range("L74").Value = 35559416
$ws.Range("M74").Value = 87.76666
$ws.Range("N74").Value = -35561164

# Sheet ARM, row 77  (@@ -11494,25 +11491,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 11853663
$ws.Range("I77").Value = 786.23334
$ws.Range("J77").Value = 35559416
$ws.Range("K77").Value = 3931.1667
$ws.Range("L77").Value = 177797080
$ws.Range("M77").Value = 436.8333000000002
$ws.Range("N77").Value = -177805816

# Sheet ARM, row 102  (@@ -12707,25 +12704,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2318.5862
$ws.Range("I102").Value = 1881.4736
$ws.Range("J102").Value = 3149.1
$ws.Range("K102").Value = 1881.4736
$ws.Range("L102").Value = 3149.1
$ws.Range("M102").Value = -259.4736
$ws.Range("N102").Value = -6393.1

# Sheet ARM, row 132  (@@ -14150,25 +14147,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 22938168
$ws.Range("I132").Value = 28385528
$ws.Range("J132").Value = 8930668
$ws.Range("K132").Value = 85156584
$ws.Range("L132").Value = 26792004
$ws.Range("M132").Value = -85154054
$ws.Range("N132").Value = -26797064

# Sheet ARM, row 134  (@@ -14251,22 +14248,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 134879.72
$ws.Range("J134").Value = 134879.72
$ws.Range("L134").Value = 134879.72
$ws.Range("N134").Value = -145019.72

# Sheet ARM, row 136  (@@ -14349,25 +14346,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2372769
$ws.Range("I136").Value = 1386163.2
$ws.Range("J136").Value = 6538438
$ws.Range("K136").Value = 4158489.6
$ws.Range("L136").Value = 19615314
$ws.Range("M136").Value = -4155939.6
$ws.Range("N136").Value = -19620414

# Sheet BSM, row 20  (@@ -15619,25 +15616,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 21752666
$ws.Range("I20").Value = 33342796
$ws.Range("J20").Value = 21172.5
$ws.Range("K20").Value = 33342796
$ws.Range("L20").Value = 21172.5
$ws.Range("M20").Value = -33342549
$ws.Range("N20").Value = -21666.5

# Sheet CRP, row 31  (@@ -23055,22 +23052,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4390982.5
$ws.Range("I31").Value = 6945806
$ws.Range("K31").Value = 6945806
$ws.Range("M31").Value = -6945511

# Sheet CRP, row 34  (@@ -23202,22 +23199,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4390982.5
$ws.Range("I34").Value = 6945806
$ws.Range("K34").Value = 6945806
$ws.Range("M34").Value = -6945604

# Sheet CRP, row 50  (@@ -23986,22 +23983,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 12977.5
$ws.Range("J50").Value = 12977.5
$ws.Range("L50").Value = 12977.5
$ws.Range("N50").Value = -14227.5

# Sheet CRP, row 51  (@@ -24035,25 +24032,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 30288.223
$ws.Range("J51").Value = 10370.571
$ws.Range("L51").Value = 10370.571
$ws.Range("N51").Value = -11842.571

# Sheet CRP, row 58  (@@ -24381,22 +24378,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1522854.6
$ws.Range("I58").Value = 11799.5
$ws.Range("J58").Value = 3789437.2
$ws.Range("K58").Value = 11799.5
$ws.Range("L58").Value = 3789437.2
$ws.Range("M58").Value = -11596.5
$ws.Range("N58").Value = -3789843.2

# Sheet CRP, row 59  (@@ -24430,22 +24430,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 17749.75
$ws.Range("J59").Value = 17749.75
$ws.Range("L59").Value = 17749.75
$ws.Range("N59").Value = -20039.75

# Sheet CRP, row 60  (@@ -24479,25 +24479,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 15739.866
$ws.Range("J60").Value = 9609.799999999999
$ws.Range("L60").Value = 9609.799999999999
$ws.Range("N60").Value = -10631.8

# Sheet CRP, row 61  (@@ -24531,25 +24531,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 30288.223
$ws.Range("J61").Value = 10370.571
$ws.Range("L61").Value = 10370.571
$ws.Range("N61").Value = -11066.571

# Sheet CRP, row 68  (@@ -24871,22 +24871,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17780.727
$ws.Range("J68").Value = 17780.727
$ws.Range("L68").Value = 17780.727
$ws.Range("N68").Value = -19278.727

# Sheet CRP, row 71  (@@ -25018,22 +25018,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 17780.727
$ws.Range("J71").Value = 17780.727
$ws.Range("L71").Value = 53342.181
$ws.Range("N71").Value = -60830.181

# Sheet CRP, row 74  (@@ -25165,25 +25165,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 17327.428
$ws.Range("I74").Value = 1285
$ws.Range("J74").Value = 18561.46
$ws.Range("K74").Value = 1285
$ws.Range("L74").Value = 18561.46
$ws.Range("M74").Value = -411
$ws.Range("N74").Value = -20309.46

# Sheet CRP, row 77  (@@ -25315,25 +25315,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 17327.428
$ws.Range("I77").Value = 1285
$ws.Range("J77").Value = 18561.46
$ws.Range("K77").Value = 3855
$ws.Range("L77").Value = 55684.38
$ws.Range("M77").Value = 513
$ws.Range("N77").Value = -64420.38

# Sheet CRP, row 103  (@@ -26592,22 +26592,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 9250
$ws.Range("I103").Value = 7666.6665
$ws.Range("J103").Value = 14000
$ws.Range("K103").Value = 7666.6665
$ws.Range("L103").Value = 14000
$ws.Range("M103").Value = -6494.6665
$ws.Range("N103").Value = -16344

# Sheet CRP, row 136  (@@ -28194,22 +28197,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1522854.6
$ws.Range("I136").Value = 11799.5
$ws.Range("J136").Value = 3789437.2
$ws.Range("K136").Value = 35398.5
$ws.Range("L136").Value = 11368311.6
$ws.Range("M136").Value = -32848.5
$ws.Range("N136").Value = -11373411.6

# Sheet GSM, row 132  (@@ -42158,25 +42164,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 12540446
$ws.Range("I132").Value = 16509874
$ws.Range("J132").Value = 7578660
$ws.Range("K132").Value = 49529622
$ws.Range("L132").Value = 22735980
$ws.Range("M132").Value = -49527092
$ws.Range("N132").Value = -22741040

# Sheet LTW, row 132  (@@ -49064,22 +49070,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5720369.5
$ws.Range("I132").Value = 7943452
$ws.Range("K132").Value = 23830356
$ws.Range("M132").Value = -23827826

# Sheet LTW, row 136  (@@ -49263,25 +49269,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10419053
$ws.Range("I136").Value = 41667332
$ws.Range("J136").Value = 2959.4443
$ws.Range("K136").Value = 125001996
$ws.Range("L136").Value = 8878.332900000001
$ws.Range("M136").Value = -124999446
$ws.Range("N136").Value = -13978.3329

# Sheet WVR, row 126  (@@ -55727,25 +55733,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1350
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1350
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 4050
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -8990

# Sheet WVR, row 136  (@@ -56223,25 +56226,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1353.3438
$ws.Range("I136").Value = 1103.4
$ws.Range("J136").Value = 5102.5
$ws.Range("K136").Value = 3310.2
$ws.Range("L136").Value = 15307.5
$ws.Range("M136").Value = -760.2000000000003
$ws.Range("N136").Value = -20407.5
